$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.016.09"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.62%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.591.30"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.61%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.480"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.50%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  -0.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.97"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.42%  "
$ws.Range("E11").Value = "  +2.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.814.11"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.70%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.593.68"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.83%  "
$ws.Range("E14").Value = "  -0.93%  "
$ws.Range("E15").Value = "  +0.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.023.70"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.80%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.17"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.85%  "
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("E19").Value = "  -0.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "201.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.63%  "
$ws.Range("E22").Value = "  -1.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +14.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.62%  "
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("E27").Value = "  -8.04%  "
$ws.Range("E28").Value = "  +0.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.48"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.68%  "
$ws.Range("E30").Value = "  +0.44%  "
$ws.Range("E31").Value = "  +0.97%  "
$ws.Range("E32").Value = "  -0.13%  "
$ws.Range("E33").Value = "  -2.86%  "
$ws.Range("E34").Value = "  -0.79%  "
$ws.Range("E35").Value = "  -0.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.129.06"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.88%  "
$ws.Range("E37").Value = "  +8.38%  "
$ws.Range("E38").Value = "  -0.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.790"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.70%  "
$ws.Range("E40").Value = "  -1.28%  "
$ws.Range("E41").Value = "  -1.94%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.780"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.85%  "
$ws.Range("E43").Value = "  -0.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.724.56"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.09"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.69%  "
$ws.Range("E46").Value = "  -0.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "53.58"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0503"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.407"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.23%  "
$ws.Range("E50").Value = "  -0.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₇0918"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -18.01%  "
